# Update "想去人数" (number of people interested) figures on the
# "展览" and "全部类型" sheets to match the newly generated data.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 5671
    $ws.Range("F3").Value = 5
    $ws.Range("F5").Value = 956
    $ws.Range("F6").Value = 29
}
